$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, pushing the existing rows 114-180 down to 115-181.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new weekly price-report entry.
# (Same market/category metadata as the surrounding rows; only the date (D) and
# volume (J) differ for this new observation.)
$ws.Cells.Item(114, 1).Value = 4
$ws.Cells.Item(114, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(114, 3).Value = "Los Lagos"
$ws.Cells.Item(114, 4).Value = 44460
$ws.Cells.Item(114, 5).Value = 10
$ws.Cells.Item(114, 6).Value = 100114014
$ws.Cells.Item(114, 7).Value = "Betarraga"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 500
$ws.Cells.Item(114, 11).Value = 1000
$ws.Cells.Item(114, 12).Value = 1000
$ws.Cells.Item(114, 13).Value = 1000
$ws.Cells.Item(114, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(114, 15).Value = "Región del Maule"
$ws.Cells.Item(114, 16).Value = 200
$ws.Cells.Item(114, 17).Value = 5
$ws.Cells.Item(114, 18).Value = "Hortaliza"
